$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "28.189.29"
$ws.Range("E2").Formula = "  +0.12%  "
$ws.Range("D3").Formula = "1.868.80"
$ws.Range("E3").Formula = "  +2.29%  "
$ws.Range("D4").Formula = "1.002"
$ws.Range("E4").Formula = "  +0.23%  "
$ws.Range("D5").Formula = "311.61"
$ws.Range("E5").Formula = "  +0.31%  "
$ws.Range("D6").Formula = "1.001"
$ws.Range("E6").Formula = "  +0.07%  "
$ws.Range("D7").Formula = "0.5047"
$ws.Range("E7").Formula = "  -1.65%  "
$ws.Range("D8").Formula = "0.3915"
$ws.Range("E8").Formula = "  -0.15%  "
$ws.Range("D9").Formula = "0.09638"
$ws.Range("E9").Formula = "  +0.82%  "
$ws.Range("D10").Formula = "1.136"
$ws.Range("E10").Formula = "  +2.62%  "
$ws.Range("D11").Formula = "40.85"
$ws.Range("E11").Formula = "  -0.17%  "
$ws.Range("D12").Formula = "6.494"
$ws.Range("E12").Formula = "  +0.58%  "
$ws.Range("D13").Formula = "20.95"
$ws.Range("E13").Formula = "  +1.97%  "
$ws.Range("D14").Formula = "1.875.17"
$ws.Range("E14").Formula = "  +3.02%  "
$ws.Range("D15").Formula = "1.002"
$ws.Range("E15").Formula = "  +0.29%  "
$ws.Range("D16").Formula = "7.416"
$ws.Range("E16").Formula = "  +0.72%  "
$ws.Range("D17").Formula = "0.00001127"
$ws.Range("E17").Formula = "  -0.55%  "
$ws.Range("D18").Formula = "92.95"
$ws.Range("E18").Formula = "  +0.37%  "
$ws.Range("D19").Formula = "0.06629"
$ws.Range("E19").Formula = "  +0.62%  "
$ws.Range("D20").Formula = "17.53"
$ws.Range("E20").Formula = "  +1.25%  "
$ws.Range("E21").Formula = "  +0.07%  "
$ws.Range("D22").Formula = "6.149"
$ws.Range("E22").Formula = "  +1.76%  "
$ws.Range("D23").Formula = "28.249.73"
$ws.Range("E23").Formula = "  +0.07%  "
$ws.Range("D24").Formula = "11.31"
$ws.Range("E24").Formula = "  +1.40%  "
$ws.Range("D25").Formula = "2.282"
$ws.Range("E25").Formula = "  +2.45%  "
$ws.Range("D26").Formula = "2.524"
$ws.Range("E26").Formula = "  +3.30%  "
$ws.Range("D27").Formula = "2.088.77"
$ws.Range("E27").Formula = "  +3.16%  "
$ws.Range("D28").Formula = "21.18"
$ws.Range("E28").Formula = "  +2.81%  "
$ws.Range("D29").Formula = "157.43"
$ws.Range("E29").Formula = "  +0.06%  "
$ws.Range("D30").Formula = "127.39"
$ws.Range("E30").Formula = "  -0.89%  "
$ws.Range("D31").Formula = "1.067"
$ws.Range("E31").Formula = "  +0.72%  "
$ws.Range("E32").Formula = "  -3.35%  "
$ws.Range("D33").Formula = "5.618"
$ws.Range("E33").Formula = "  -0.55%  "
$ws.Range("D34").Formula = "3.629"
$ws.Range("E34").Formula = "  -0.04%  "
$ws.Range("D35").Formula = "9.566"
$ws.Range("E35").Formula = "  +5.03%  "
$ws.Range("D36").Formula = "0.06748"
$ws.Range("E36").Formula = "  -2.56%  "
$ws.Range("D37").Formula = "0.02384"
$ws.Range("E37").Formula = "  +1.75%  "
$ws.Range("D38").Formula = "0.2179"
$ws.Range("E38").Formula = "  +0.08%  "
$ws.Range("D39").Formula = "0.6349"
$ws.Range("E39").Formula = "  +2.30%  "
$ws.Range("D40").Formula = "11.45"
$ws.Range("E40").Formula = "  -1.05%  "
$ws.Range("D41").Formula = "4.984"
$ws.Range("E41").Formula = "  -0.84%  "
$ws.Range("D42").Formula = "1.181"
$ws.Range("E42").Formula = "  +2.50%  "
$ws.Range("D43").Formula = "1.001"
$ws.Range("E43").Formula = "  +0.18%  "
$ws.Range("D44").Formula = "13.60"
$ws.Range("E44").Formula = "  +1.94%  "
$ws.Range("D45").Formula = "0.6014"
$ws.Range("E45").Formula = "  +0.70%  "
$ws.Range("D46").Formula = "3.661"
$ws.Range("E46").Formula = "  -1.25%  "
$ws.Range("E47").Formula = "  -2.02%  "
$ws.Range("D48").Formula = "124.14"
$ws.Range("E48").Formula = "  -0.89%  "
$ws.Range("D49").Formula = "1.992"
$ws.Range("E49").Formula = "  +1.16%  "
$ws.Range("D50").Formula = "1.194"
$ws.Range("E50").Formula = "  +0.66%  "
$ws.Range("D51").Formula = "0.06835"
$ws.Range("E51").Formula = "  +0.82%  "
